# Daily attendance processing - swap order of "System" entry with the
# adjacent recorder email in the "Recorded By" column (G).
#
# Rule (derived from the target diff): for every cell in column G whose
# value is a comma-separated list ("a, b" or "a, b, c") that contains the
# exact token "System" as the first or second token, swap the first two
# tokens and keep any remaining tokens (e.g. a trailing lowercase
# "system") in place. Cells with a single token, or without "System" in
# the first two tokens, are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }

    $tokens = $val -split ", "

    if ($tokens.Count -ge 2 -and ($tokens[0] -eq "System" -or $tokens[1] -eq "System")) {
        $newTokens = @($tokens[1], $tokens[0])
        if ($tokens.Count -gt 2) {
            $newTokens = $newTokens + $tokens[2..($tokens.Count - 1)]
        }
        $cell.Value2 = [string]::Join(", ", $newTokens)
    }
}
